# Update "想去人数" (want-to-go count) values in column F across the
# workbook's four sheets, as generated by the gh-pages build at 456a3b4.

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 2927
$ws1.Range("F6").Value = 2927
$ws1.Range("F7").Value = 788
$ws1.Range("F10").Value = 412
$ws1.Range("F13").Value = 502
$ws1.Range("F15").Value = 2153
$ws1.Range("F26").Value = 595
$ws1.Range("F27").Value = 595
$ws1.Range("F28").Value = 9
$ws1.Range("F32").Value = 565

# 演出 (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F27").Value = 104
$ws2.Range("F36").Value = 529
$ws2.Range("F37").Value = 529
$ws2.Range("F38").Value = 17

# 本地生活 (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 65
$ws3.Range("F5").Value = 569
$ws3.Range("F6").Value = 230
$ws3.Range("F7").Value = 247

# 全部类型 (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 569
$ws4.Range("F6").Value = 230
$ws4.Range("F10").Value = 2927
$ws4.Range("F11").Value = 2927
$ws4.Range("F12").Value = 788
$ws4.Range("F15").Value = 412
$ws4.Range("F19").Value = 502
$ws4.Range("F25").Value = 2153
$ws4.Range("F33").Value = 247
$ws4.Range("F38").Value = 595
$ws4.Range("F39").Value = 595
$ws4.Range("F42").Value = 565
$ws4.Range("F50").Value = 529
$ws4.Range("F51").Value = 17
